# Commit: "Fixed interface to show Friesland, added the whole of the Netherlands"
#
# Data change: insert a new row 2 in the "area" lookup sheet that maps
# "Nederland" -> "Nederland" (so the whole country can be selected as an
# aggregate "area", same way a municipality maps to its province), pushing
# every existing municipality/province row down by one.

$wb = $excel.ActiveWorkbook

$wsOutcome = $wb.Worksheets.Item("outcome")
$wsArea    = $wb.Worksheets.Item("area")

# Insert a fresh row above the first data row (row 2, right under the header).
$wsArea.Rows.Item(2).Insert()

# Fill in the new row. Column A gets plain (unstyled) text, column B is
# styled like the header row (bold-ish "type" font), matching how the
# existing province names are emphasized elsewhere in the sheet.
$wsArea.Range("A2").Value = "Nederland"
$wsArea.Range("B2").Value = "Nederland"
$wsArea.Range("A2").ClearFormats()

# --- View / interface state -------------------------------------------------
# The workbook now opens on the "area" sheet (interface shows Friesland /
# the rest of the area list instead of being stuck scrolled into "outcome").
$wsOutcome.Activate()
$wsOutcome.Range("A2").Select()

$wsArea.Activate()
$wsArea.Range("B2").Select()
